# The header ("EMPRESA: MKTO empresa IMPORTACIONES, S.L.") is corrected
# to just show the "empresa" bookmark field, removing the hard-coded
# "EMPRESA:", "MKTO" and " IMPORTACIONES, S.L." literal text runs while
# keeping the bookmark ("empresa") and the surrounding table intact.
#
# Targeted Find/Replace on a sub-range of a header does not reliably
# persist in this runtime, so we rebuild the header's primary story by
# replacing its whole Range with equivalent WordprocessingML via
# InsertXML, which *does* persist and lets us control the exact
# resulting markup (bookmark, run formatting, table, trailing
# paragraph, etc.) precisely.

$d = $word.ActiveDocument

$targetHeader = $null
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Range.Text -like "*MKTO*") {
        $targetHeader = $hdr
    }
}

if ($targetHeader -ne $null) {
    $rng = $targetHeader.Range
    $xml = @'
<w:tbl><w:tblPr><w:tblW w:w="9540" w:type="dxa"/><w:jc w:val="center"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblLayout w:type="fixed"/><w:tblCellMar><w:left w:w="70" w:type="dxa"/><w:right w:w="70" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="3240"/><w:gridCol w:w="6300"/></w:tblGrid><w:tr w:rsidR="00D123F0" w14:paraId="498A2A05" w14:textId="77777777" w:rsidTr="009F5449"><w:trPr><w:cantSplit/><w:trHeight w:val="1105"/><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3240" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="444CAC6A" w14:textId="77777777" w:rsidR="00D123F0" w:rsidRDefault="00D123F0" w:rsidP="00984760"><w:pPr><w:pStyle w:val="Encabezado"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr><w:t>INFORME DE INVESTIGACIÓN DE ACCIDENTES / INCIDENTES.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6300" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="4C7F4516" w14:textId="387A7591" w:rsidR="00D123F0" w:rsidRPr="009F5449" w:rsidRDefault="00D123F0" w:rsidP="00926A69"><w:pPr><w:pStyle w:val="Encabezado"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:bookmarkStart w:id="14" w:name="empresa"/><w:r w:rsidR="00476151"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:szCs w:val="28"/></w:rPr><w:t>empresa</w:t></w:r><w:bookmarkEnd w:id="14"/></w:p></w:tc></w:tr></w:tbl>
'@
    [void]$rng.InsertXML($xml)
}
